$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 70 is no longer the last row, so it gets the regular date/time number
# format (same style previously used by row 69) instead of the special
# "last row" date-only format.
$ws.Range("A70").NumberFormat = $ws.Range("A69").NumberFormat

# Add the new daily row (71) with the "last row" date-only style/format
# that row 70 used to have.
$ws.Range("A71").NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(71, 1).Value = 45658
$ws.Cells.Item(71, 2).Value = 166
$ws.Cells.Item(71, 3).Value = 162
$ws.Cells.Item(71, 4).Value = 164
